# QC code and plot updates
#
# 1) treatments.csv (sheet1): insert two new "Ambient" rows at the top of the
#    data block (row 3/4) for the Set=4 and Set=22 ambient sensors, pushing
#    the existing 30 treatment rows down by two.
# 2) valvemap.csv (sheet2): the generic "Ambient" treatment label used for
#    MPVPosition 10 rows is disambiguated to "Ambient4" / "Ambient22"
#    depending on which physical loop (Set) the reading came from.
# 3) Cosmetic: window size/position + active-cell selection on both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "treatments.csv"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("treatments.csv")

# Insert two blank rows above the first data row (row 3), shifting rows
# 3-32 down to 5-34.
$ws1.Rows.Item(3).Insert()
$ws1.Rows.Item(3).Insert()

# New row 3: Ambient / Field moisture-style "Ambient" category, Set 4
$ws1.Range("A3").Value = "Ambient4"
$ws1.Range("B3").Value = "Ambient"
$ws1.Range("C3").Value = 4

# New row 4: same, Set 22
$ws1.Range("A4").Value = "Ambient22"
$ws1.Range("B4").Value = "Ambient"
$ws1.Range("C4").Value = 22

# Restore the active-cell selection on this sheet.
$ws1.Activate()
[void]$ws1.Range("A4").Select()

# ---------------------------------------------------------------------
# Sheet "valvemap.csv"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("valvemap.csv")

# Rows whose E column says the generic "Ambient" treatment, now split by
# which physical loop (Set, column C) produced the reading.
$ambient22Rows = @(3, 10, 17, 38, 45, 52, 73, 80, 87, 122, 129, 136)
$ambient4Rows  = @(24, 31, 59, 66, 94, 101, 108, 115)

foreach ($r in $ambient22Rows) {
    $ws2.Range("E$r").Value = "Ambient22"
}
foreach ($r in $ambient4Rows) {
    $ws2.Range("E$r").Value = "Ambient4"
}

# Restore selection on this sheet.
$ws2.Activate()
[void]$ws2.Range("E1").Select()

# ---------------------------------------------------------------------
# Cosmetic window geometry (workbook-level view).
# ---------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 7540
$win.Top = 0
$win.Width = 24500
$win.Height = 15360

$ws2.Activate()
